# Auto-generated Excel COM-interop script to apply scheduled runner updates
# to the Bahamut_Profits Leve-profit tables across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 23257990
$ws.Range("J40").Value = 33335636
$ws.Range("L40").Value = 33335636
$ws.Range("N40").Value = -33335986

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4644.4443
$ws.Range("I64").Value = 5400
$ws.Range("J64").Value = 4040
$ws.Range("K64").Value = 5400
$ws.Range("L64").Value = 4040
$ws.Range("M64").Value = -5152
$ws.Range("N64").Value = -4536

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4644.4443
$ws.Range("I67").Value = 5400
$ws.Range("J67").Value = 4040
$ws.Range("K67").Value = 5400
$ws.Range("L67").Value = 4040
$ws.Range("M67").Value = -4542
$ws.Range("N67").Value = -5756

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 939.6316
$ws.Range("I137").Value = 881.25
$ws.Range("J137").Value = 1251
$ws.Range("K137").Value = 2643.75
$ws.Range("L137").Value = 3753
$ws.Range("M137").Value = -93.75
$ws.Range("N137").Value = -8853

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 8025
$ws.Range("I34").Value = 8025
$ws.Range("K34").Value = 8025
$ws.Range("M34").Value = -7754

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1654.25
$ws.Range("I61").Value = 1593.7
$ws.Range("J61").Value = 1957
$ws.Range("K61").Value = 1593.7
$ws.Range("L61").Value = 1957
$ws.Range("M61").Value = -1381.7
$ws.Range("N61").Value = -2381

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1120.4333
$ws.Range("I74").Value = 1154.25
$ws.Range("J74").Value = 1052.8
$ws.Range("K74").Value = 1154.25
$ws.Range("L74").Value = 1052.8
$ws.Range("M74").Value = -280.25
$ws.Range("N74").Value = -2800.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1120.4333
$ws.Range("I77").Value = 1154.25
$ws.Range("J77").Value = 1052.8
$ws.Range("K77").Value = 5771.25
$ws.Range("L77").Value = 5264
$ws.Range("M77").Value = -1403.25
$ws.Range("N77").Value = -14000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1452.5
$ws.Range("I132").Value = 983.16
$ws.Range("K132").Value = 2949.48
$ws.Range("M132").Value = -419.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1654.25
$ws.Range("I136").Value = 1593.7
$ws.Range("J136").Value = 1957
$ws.Range("K136").Value = 4781.1
$ws.Range("L136").Value = 5871
$ws.Range("M136").Value = -2231.1
$ws.Range("N136").Value = -10971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 350.2
$ws.Range("I5").Value = 78
$ws.Range("J5").Value = 531.6667
$ws.Range("K5").Value = 78
$ws.Range("L5").Value = 531.6667
$ws.Range("M5").Value = 35
$ws.Range("N5").Value = -757.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8051.1
$ws.Range("I105").Value = 8600
$ws.Range("J105").Value = 7227.75
$ws.Range("K105").Value = 8600
$ws.Range("L105").Value = 7227.75
$ws.Range("M105").Value = -6853
$ws.Range("N105").Value = -10721.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27435.666
$ws.Range("I134").Value = 4940.0586
$ws.Range("J134").Value = 42732.68
$ws.Range("K134").Value = 14820.1758
$ws.Range("L134").Value = 128198.04
$ws.Range("M134").Value = -12285.1758
$ws.Range("N134").Value = -133268.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 7
$ws.Range("K10").Value = 7
$ws.Range("M10").Value = 132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2513.7307
$ws.Range("I31").Value = 2462.28
$ws.Range("J31").Value = 3800
$ws.Range("K31").Value = 2462.28
$ws.Range("L31").Value = 3800
$ws.Range("M31").Value = -2167.28
$ws.Range("N31").Value = -4390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2513.7307
$ws.Range("I34").Value = 2462.28
$ws.Range("J34").Value = 3800
$ws.Range("K34").Value = 2462.28
$ws.Range("L34").Value = 3800
$ws.Range("M34").Value = -2260.28
$ws.Range("N34").Value = -4204

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6421.0527
$ws.Range("I58").Value = 671.25
$ws.Range("J58").Value = 10602.728
$ws.Range("K58").Value = 671.25
$ws.Range("L58").Value = 10602.728
$ws.Range("M58").Value = -468.25
$ws.Range("N58").Value = -11008.728

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2359.3076
$ws.Range("I132").Value = 1823.6666
$ws.Range("J132").Value = 4144.778
$ws.Range("K132").Value = 5470.9998
$ws.Range("L132").Value = 12434.334
$ws.Range("M132").Value = -2940.9998
$ws.Range("N132").Value = -17494.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1918.6976
$ws.Range("I134").Value = 1656.6552
$ws.Range("J134").Value = 2461.5
$ws.Range("K134").Value = 4969.9656
$ws.Range("L134").Value = 7384.5
$ws.Range("M134").Value = -2434.9656
$ws.Range("N134").Value = -12454.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6421.0527
$ws.Range("I136").Value = 671.25
$ws.Range("J136").Value = 10602.728
$ws.Range("K136").Value = 2013.75
$ws.Range("L136").Value = 31808.184
$ws.Range("M136").Value = 536.25
$ws.Range("N136").Value = -36908.18399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 73.333336
$ws.Range("I40").Value = 73.333336
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 293.333344
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -224.333344
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 939.75
$ws.Range("I133").Value = 939.75
$ws.Range("K133").Value = 2819.25
$ws.Range("M133").Value = 2240.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1662.3684
$ws.Range("I132").Value = 753.36365
$ws.Range("J132").Value = 2912.25
$ws.Range("K132").Value = 2260.09095
$ws.Range("L132").Value = 8736.75
$ws.Range("M132").Value = 269.9090500000002
$ws.Range("N132").Value = -13796.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2373.5925
$ws.Range("I132").Value = 1737.6
$ws.Range("J132").Value = 3545.158
$ws.Range("K132").Value = 5212.799999999999
$ws.Range("L132").Value = 10635.474
$ws.Range("M132").Value = -2682.799999999999
$ws.Range("N132").Value = -15695.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3844
$ws.Range("I136").Value = 1214.3182
$ws.Range("J136").Value = 11075.625
$ws.Range("K136").Value = 3642.9546
$ws.Range("L136").Value = 33226.875
$ws.Range("M136").Value = -1092.9546
$ws.Range("N136").Value = -38326.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2093.5173
$ws.Range("I132").Value = 1652.8096
$ws.Range("J132").Value = 3250.375
$ws.Range("K132").Value = 4958.4288
$ws.Range("L132").Value = 9751.125
$ws.Range("M132").Value = -2428.4288
$ws.Range("N132").Value = -14811.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1908.92
$ws.Range("I136").Value = 1827.0869
$ws.Range("J136").Value = 2850
$ws.Range("K136").Value = 5481.2607
$ws.Range("L136").Value = 8550
$ws.Range("M136").Value = -2931.2607
$ws.Range("N136").Value = -13650

Write-Host "Applied Bahamut_Profits scheduled update across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
